$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LOM3219")

# Créditos-aula: 4 -> 2 (keep stored as text, like original)
$ws.Range("B5").NumberFormat = "@"
$ws.Range("C5").NumberFormat = "@"
$ws.Range("B5").Value = "2"
$ws.Range("C5").Value = "2"

# Carga horária: 60 h -> 30 h
$ws.Range("B7").Value = "30 h"
$ws.Range("C7").Value = "30 h"

# Ativação: 01/01/2012 -> 01/01/2023 (keep stored as text, like original)
$ws.Range("B8").NumberFormat = "@"
$ws.Range("C8").NumberFormat = "@"
$ws.Range("B8").Value = "01/01/2023"
$ws.Range("C8").Value = "01/01/2023"

# Objetivos: responsible professor changed
$ws.Range("B10").Value = "7290967 - Emerson Gonçalves de Melo"
$ws.Range("C10").Value = "7290967 - Emerson Gonçalves de Melo"

# Objectives: add English objectives text (row 11, was empty on B/C)
$ws.Range("B11").Value = "Present the concepts of nanoscience and nanotechnology. The physical and chemical properties of materials on a nanometer scale are described by the laws of quantum mechanics, presenting in these dimensions different characteristics of materials on a macroscopic scale. The knowledge of this interdisciplinary area is fundamental in the formation of a researcher and/or a professional working in the area of materials."
$ws.Range("C11").Value = "Present the concepts of nanoscience and nanotechnology. The physical and chemical properties of materials on a nanometer scale are described by the laws of quantum mechanics, presenting in these dimensions different characteristics of materials on a macroscopic scale. The knowledge of this interdisciplinary area is fundamental in the formation of a researcher and/or a professional working in the area of materials."
$ws.Cells.Item(11, 2).Style = $ws.Cells.Item(10, 2).Style
$ws.Cells.Item(11, 3).Style = $ws.Cells.Item(10, 3).Style

# Programa resumido: "Semestral" -> "01/01/2023" (reuses Ativação string, stays text)
$ws.Range("B13").NumberFormat = "@"
$ws.Range("C13").NumberFormat = "@"
$ws.Range("B13").Value = "01/01/2023"
$ws.Range("C13").Value = "01/01/2023"

# Short syllabus: add short syllabus text (row 14, was empty on B/C)
$ws.Range("B14").Value = "Nanoscience and nanotechnology: principles and applications."
$ws.Range("C14").Value = "Nanoscience and nanotechnology: principles and applications."
$ws.Cells.Item(14, 2).Style = $ws.Cells.Item(13, 2).Style
$ws.Cells.Item(14, 3).Style = $ws.Cells.Item(13, 3).Style

# Programa: "01/01/2012" -> "7290967 - Emerson Gonçalves de Melo" (reuses Objetivos string)
$ws.Range("B15").Value = "7290967 - Emerson Gonçalves de Melo"
$ws.Range("C15").Value = "7290967 - Emerson Gonçalves de Melo"

# Syllabus: add full syllabus text (row 16, was empty on B/C)
$ws.Range("B16").Value = "Conceptualization: nanoscience and nanotechnology. Low-dimensional systems. Quantum Confinement. Chemical bonds: molecules and clusters. Electronic and structural properties. Synthesis and fabrication of materials at the nanometer scale: bottom-up and top-down techniques. Fullerenes and carbon nanotubes. Molecular self-organization and supramolecular systems. Quantum wires and dots. Magnetic nanoparticles. Characterization techniques: X-ray diffraction, scattering and absorption, scanning tunneling microscopy (STM), atomic force microscopy (AFM), transmission electron microscopy. Transport properties: ballistic transport, quantum conductance, Coulomb blocking. Molecular devices. Diffusive transport. Nanomagnetism: magnetic order, superparamagnetism and spintronics. Applications."
$ws.Range("C16").Value = "Conceptualization: nanoscience and nanotechnology. Low-dimensional systems. Quantum Confinement. Chemical bonds: molecules and clusters. Electronic and structural properties. Synthesis and fabrication of materials at the nanometer scale: bottom-up and top-down techniques. Fullerenes and carbon nanotubes. Molecular self-organization and supramolecular systems. Quantum wires and dots. Magnetic nanoparticles. Characterization techniques: X-ray diffraction, scattering and absorption, scanning tunneling microscopy (STM), atomic force microscopy (AFM), transmission electron microscopy. Transport properties: ballistic transport, quantum conductance, Coulomb blocking. Molecular devices. Diffusive transport. Nanomagnetism: magnetic order, superparamagnetism and spintronics. Applications."
$ws.Cells.Item(16, 2).Style = $ws.Cells.Item(15, 2).Style
$ws.Cells.Item(16, 3).Style = $ws.Cells.Item(15, 3).Style

# Método: 519033 - Carlos Yujiro Shigue -> 1176388 - Luiz Tadeu Fernandes Eleno
$ws.Range("B18").Value = "1176388 - Luiz Tadeu Fernandes Eleno"
$ws.Range("C18").Value = "1176388 - Luiz Tadeu Fernandes Eleno"
